# Apply weekly fruit/vegetable price data update
# (rows 2-13, columns D, M, N, O, P, S are updated to reflect the latest
#  weekly reshuffle of dates / volumes / prices)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: Date(D), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), PrecioKg(S)
$data = @{
    2  = @{ D = 44914; M = 56; N = 23000; O = 23000; P = 23000; S = 4600 }
    3  = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    4  = @{ D = 44907; M = 45; N = 25000; O = 25000; P = 25000; S = 5000 }
    5  = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    6  = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    7  = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    8  = @{ D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 }
    9  = @{ D = 44931; M = 50; N = 18000; O = 18000; P = 18000; S = 3600 }
    10 = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    11 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    12 = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
    13 = @{ D = 44902; M = 35; N = 12000; O = 12000; P = 12000; S = 2400 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D
    $ws.Cells.Item($row, 13).Value = $vals.M   # M
    $ws.Cells.Item($row, 14).Value = $vals.N   # N
    $ws.Cells.Item($row, 15).Value = $vals.O   # O
    $ws.Cells.Item($row, 16).Value = $vals.P   # P
    $ws.Cells.Item($row, 19).Value = $vals.S   # S
}
